$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mandreoli")

# Column A (Datum) - each date moves down one row, with a new date appended for row 5
$ws.Range("A3").Value = "16.01.2019"
$ws.Range("A4").Value = "20.01.2019"
$ws.Range("A5").Value = "25.01.2019"

# Column D (Zeit) holds text-formatted numbers ("0.50", "1.00", ...). Force text
# formatting first so the values stay text (not auto-converted to numbers),
# matching the original "1.50"/"0.50"/"1.00" style of entry.
$ws.Range("D3:D5").NumberFormat = "@"
$ws.Range("D3").Value = "0.50"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "2.00"

# Column E (Zuschlag) holds plain numeric surcharge multipliers
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1.25
$ws.Range("E5").Value = 1
